$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for account 000772433 / MARCELO / 5000 (Excel row 5)
$ws.Rows.Item(5).Delete()

# Insert a new row after account 004220849 / DULCE (now Excel row 14, after the
# deletion above) and before 004432579 / ANA (now Excel row 15).
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "008002502"
$ws.Cells.Item(15, 2).Value = "JORGEANA"
$ws.Cells.Item(15, 3).Value = 500
